# Generate Report for Handback
# Adds a new handback record (6eb2f1be-a51f-4f10-b6a6-1c6633556a4d) as row 3
# to the "Overview", "zh-cn" and "de-de" sheets/tables.

$wb = $excel.ActiveWorkbook

$newFile   = "6eb2f1be-a51f-4f10-b6a6-1c6633556a4d.md"
$newPath   = "e2e\6eb2f1be-a51f-4f10-b6a6-1c6633556a4d.md"
$status    = "Handed back: in sync with en-US"
$ext       = ".md"
$genDate   = "2016-09-06 14:05:43"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Set-TextValue {
    param($range, [string]$text)
    # Force text storage so look-alike booleans ("True"/"False") and other
    # auto-detected literals are not silently converted to other types.
    $range.NumberFormat = "@"
    $range.Value2 = $text
}

function Set-DateTextValue {
    param($range, [string]$text)
    # Keep the value as plain text (matches source data) but apply the same
    # display format used by the other date/time columns.
    $range.NumberFormat = $dateFmt
    $range.Value2 = $text
}

# ---------------------------------------------------------------------
# Sheet "Overview"  (columns: A,B,C,D,E,F,G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-TextValue $wsOverview.Range("A3") $newFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a56e3e229e61049f3e3a04c46fa853cb05b3b844/e2e/$newFile", "", "", $newPath)
Set-TextValue $wsOverview.Range("C3") $ext
Set-TextValue $wsOverview.Range("E3") $status
Set-TextValue $wsOverview.Range("F3") $status
Set-DateTextValue $wsOverview.Range("G3") $genDate

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"  (columns: A..P)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf   = "6eb2f1be-a51f-4f10-b6a6-1c6633556a4d.b986a7f453da8fccf5767a4102b4a1698c19379c.zh-cn.xlf"
$zhHoDt  = "2016-09-06 14:05:26"
$zhHbDt  = "2016-09-06 14:06:37"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a56e3e229e61049f3e3a04c46fa853cb05b3b844/e2e/$newFile", "", "", $newFile)
Set-TextValue $wsZhCn.Range("B3") $ext
Set-TextValue $wsZhCn.Range("C3") $status
Set-TextValue $wsZhCn.Range("D3") "e2e"
Set-TextValue $wsZhCn.Range("E3") "ht"
Set-TextValue $wsZhCn.Range("F3") "True"
Set-TextValue $wsZhCn.Range("G3") $zhXlf
Set-DateTextValue $wsZhCn.Range("H3") $zhHoDt
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f56e6147307652219a6dda66a45fc5e3f21c01d6/e2e/$newFile", "", "", $newFile)
Set-TextValue $wsZhCn.Range("J3") $zhXlf
Set-DateTextValue $wsZhCn.Range("K3") $zhHbDt
Set-TextValue $wsZhCn.Range("M3") "True"
Set-TextValue $wsZhCn.Range("O3") "False"

$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"  (columns: A..P)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf   = "6eb2f1be-a51f-4f10-b6a6-1c6633556a4d.b986a7f453da8fccf5767a4102b4a1698c19379c.de-de.xlf"
$deHoDt  = "2016-09-06 14:05:43"
$deHbDt  = "2016-09-06 14:06:58"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a56e3e229e61049f3e3a04c46fa853cb05b3b844/e2e/$newFile", "", "", $newFile)
Set-TextValue $wsDeDe.Range("B3") $ext
Set-TextValue $wsDeDe.Range("C3") $status
Set-TextValue $wsDeDe.Range("D3") "e2e"
Set-TextValue $wsDeDe.Range("E3") "ht"
Set-TextValue $wsDeDe.Range("F3") "True"
Set-TextValue $wsDeDe.Range("G3") $deXlf
Set-DateTextValue $wsDeDe.Range("H3") $deHoDt
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/16c3737d140a93860deefb3d80e934bbf88c9312/e2e/$newFile", "", "", $newFile)
Set-TextValue $wsDeDe.Range("J3") $deXlf
Set-DateTextValue $wsDeDe.Range("K3") $deHbDt
Set-TextValue $wsDeDe.Range("M3") "True"
Set-TextValue $wsDeDe.Range("O3") "False"

$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
